# Updates Price (D) and Volume(1h) (E) columns for the crypto symbol list,
# matching the refreshed scrape values from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original cells' text type) instead of inferring numbers/percentages,
# then formats are cleared so no stray number-format/style is left behind.
$updates = [ordered]@{
    "D2" = "258.94"
    "E2" = "0.63%"
    "D3" = "26.94"
    "E3" = "-0.03%"
    "D4" = "4.678"
    "E4" = "0.13%"
    "D5" = "0.06036"
    "E5" = "2.63%"
    "D6" = "6.680"
    "E6" = "0.79%"
    "D7" = "0.8599"
    "D8" = "0.9169"
    "E8" = "-3.41%"
    "D9" = "0.1394"
    "E9" = "-1.06%"
    "D10" = "0.05361"
    "E10" = "35.23%"
    "D11" = "0.07073"
    "E11" = "-0.31%"
    "D12" = "0.03076"
    "E12" = "-3.20%"
    "D13" = "0.09123"
    "E13" = "-0.42%"
    "D14" = "0.001523"
    "E14" = "-1.43%"
    "D15" = "0.0006042"
    "E15" = "0.23%"
    "D16" = "0.006193"
    "E16" = "-0.01%"
    "D17" = "3.465"
    "E17" = "-1.72%"
    "D18" = "3.176"
    "E18" = "-1.01%"
    "D20" = "0.3129"
    "E20" = "1.54%"
    "D21" = "0.1297"
    "D22" = "4.140"
    "E22" = "8.39%"
    "D23" = "0.04231"
    "E23" = "0.16%"
    "E24" = "-0.45%"
    "D25" = "0.004028"
    "E25" = "-6.40%"
    "D26" = "0.0001199"
    "E26" = "-0.03%"
    "E27" = "-21.33%"
    "D40" = "0.03848"
    "E40" = "0.50%"
    "D41" = "0.1114"
    "E41" = "0.96%"
    "D42" = "0.004025"
    "E42" = "2.50%"
    "D43" = "0.01530"
    "E43" = "33.31%"
    "D44" = "0.002199"
    "E44" = "-9.50%"
    "D45" = "0.00005106"
    "E45" = "-6.42%"
    "E46" = "-0.03%"
    "D47" = "0.05452"
    "E47" = "-22.10%"
    "D48" = "0.1353"
    "E48" = "-27.47%"
    "D49" = "0.00002099"
    "E49" = "-0.03%"
    "D50" = "0.0001999"
    "E50" = "-0.03%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.ClearFormats()
}
